# Update TPM-derived values for the Mfge8-Itgav ligand-receptor pair table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 27.852944
$ws.Range("H2").Value = 83.558832
$ws.Range("I2").Value = 0.2559209115167818
$ws.Range("J2").Value = 0.2559209115167818
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 245.6809962190827
$ws.Range("R2").Value = 2211.128965971744
$ws.Range("S2").Value = 0.01641770180846443
$ws.Range("T2").Value = 0.01641770180846443

# Row 3
$ws.Range("G3").Value = 27.852944
$ws.Range("H3").Value = 83.558832
$ws.Range("I3").Value = 0.2559209115167818
$ws.Range("J3").Value = 0.2559209115167818
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 1523.87732050904
$ws.Range("R3").Value = 13714.89588458136
$ws.Range("S3").Value = 0.101833531391615
$ws.Range("T3").Value = 0.101833531391615

# Row 4
$ws.Range("G4").Value = 27.852944
$ws.Range("H4").Value = 83.558832
$ws.Range("I4").Value = 0.2559209115167818
$ws.Range("J4").Value = 0.2559209115167818
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 610.2068464661867
$ws.Range("R4").Value = 5491.86161819568
$ws.Range("S4").Value = 0.04077724447938864
$ws.Range("T4").Value = 0.04077724447938864

# Row 5
$ws.Range("G5").Value = 27.852944
$ws.Range("H5").Value = 83.558832
$ws.Range("I5").Value = 0.2559209115167818
$ws.Range("J5").Value = 0.2559209115167818
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 1449.936778542896
$ws.Range("R5").Value = 13049.43100688606
$ws.Range("S5").Value = 0.0968924338373137
$ws.Range("T5").Value = 0.09689243383731372

# Row 6
$ws.Range("I6").Value = 0.3112048767201538
$ws.Range("J6").Value = 0.3112048767201538
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 298.7529377247885
$ws.Range("R6").Value = 2688.776439523096
$ws.Range("S6").Value = 0.01996424925595181
$ws.Range("T6").Value = 0.01996424925595181

# Row 7
$ws.Range("I7").Value = 0.3112048767201538
$ws.Range("J7").Value = 0.3112048767201538
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.1238315829483412
$ws.Range("T7").Value = 0.1238315829483412

# Row 8
$ws.Range("I8").Value = 0.3112048767201538
$ws.Range("J8").Value = 0.3112048767201538
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 742.0235622904578
$ws.Range("R8").Value = 6678.212060614121
$ws.Range("S8").Value = 0.04958593366202346
$ws.Range("T8").Value = 0.04958593366202346

# Row 9
$ws.Range("I9").Value = 0.3112048767201538
$ws.Range("J9").Value = 0.3112048767201538
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 1763.151724273497
$ws.Range("R9").Value = 15868.36551846147
$ws.Range("S9").Value = 0.1178231108538373
$ws.Range("T9").Value = 0.1178231108538374

# Row 10
$ws.Range("G10").Value = 30.14135433333333
$ws.Range("H10").Value = 90.42406299999999
$ws.Range("I10").Value = 0.2769474880406526
$ws.Range("J10").Value = 0.2769474880406526
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 265.8662567233717
$ws.Range("R10").Value = 2392.796310510346
$ws.Range("S10").Value = 0.01776658753013447
$ws.Range("T10").Value = 0.01776658753013447

# Row 11
$ws.Range("G11").Value = 30.14135433333333
$ws.Range("H11").Value = 90.42406299999999
$ws.Range("I11").Value = 0.2769474880406526
$ws.Range("J11").Value = 0.2769474880406526
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 1649.079762555568
$ws.Range("R11").Value = 14841.71786300011
$ws.Range("S11").Value = 0.1102002198650631
$ws.Range("T11").Value = 0.1102002198650631

# Row 12
$ws.Range("G12").Value = 30.14135433333333
$ws.Range("H12").Value = 90.42406299999999
$ws.Range("I12").Value = 0.2769474880406526
$ws.Range("J12").Value = 0.2769474880406526
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 660.341713822541
$ws.Range("R12").Value = 5943.075424402869
$ws.Range("S12").Value = 0.04412752111913963
$ws.Range("T12").Value = 0.04412752111913962

# Row 13
$ws.Range("G13").Value = 30.14135433333333
$ws.Range("H13").Value = 90.42406299999999
$ws.Range("I13").Value = 0.2769474880406526
$ws.Range("J13").Value = 0.2769474880406526
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 1569.064232599372
$ws.Range("R13").Value = 14121.57809339435
$ws.Range("S13").Value = 0.1048531595263154
$ws.Range("T13").Value = 0.1048531595263154

# Row 14
$ws.Range("G14").Value = 16.970158
$ws.Range("H14").Value = 50.910474
$ws.Range("I14").Value = 0.1559267237224118
$ws.Range("J14").Value = 0.1559267237224118
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 149.6877789089453
$ws.Range("R14").Value = 1347.190010180508
$ws.Range("S14").Value = 0.01000292800956793
$ws.Range("T14").Value = 0.01000292800956793

# Row 15
$ws.Range("G15").Value = 16.970158
$ws.Range("H15").Value = 50.910474
$ws.Range("I15").Value = 0.1559267237224118
$ws.Range("J15").Value = 0.1559267237224118
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 928.46339337253
$ws.Range("R15").Value = 8356.17054035277
$ws.Range("S15").Value = 0.0620448279152705
$ws.Range("T15").Value = 0.0620448279152705

# Row 16
$ws.Range("G16").Value = 16.970158
$ws.Range("H16").Value = 50.910474
$ws.Range("I16").Value = 0.1559267237224118
$ws.Range("J16").Value = 0.1559267237224118
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 371.7849932564734
$ws.Range("R16").Value = 3346.06493930826
$ws.Range("S16").Value = 0.02484463694824695
$ws.Range("T16").Value = 0.02484463694824694

# Row 17
$ws.Range("G17").Value = 16.970158
$ws.Range("H17").Value = 50.910474
$ws.Range("I17").Value = 0.1559267237224118
$ws.Range("J17").Value = 0.1559267237224118
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 883.413122213722
$ws.Range("R17").Value = 7950.718099923497
$ws.Range("S17").Value = 0.05903433084932638
$ws.Range("T17").Value = 0.05903433084932638

